$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric remain stored as text,
# matching the original inline-string cell type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.416.23"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.247.89"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.82"
$ws.Range("E5").Value = "  +1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.01"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.24"
$ws.Range("E10").Value = "  +3.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.25"
$ws.Range("E12").Value = "  +2.47%  "

$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.378.29"
$ws.Range("E14").Value = "  +4.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.843"
$ws.Range("E15").Value = "  +4.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.75"
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "44.137.98"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0968"
$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.40"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("E20").Value = "  +4.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.98"
$ws.Range("E21").Value = "  +2.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.04"
$ws.Range("E22").Value = "  +5.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.08"
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("E24").Value = "  +5.74%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +5.95%  "

$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.61"
$ws.Range("E27").Value = "  +8.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.90"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.00"
$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.14"
$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.48"
$ws.Range("E31").Value = "  +1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0803"
$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  -4.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +4.82%  "

$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  +5.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +8.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.76"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.85"
$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("E41").Value = "  +2.45%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.750.72"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.195"
$ws.Range("E44").Value = "  +5.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "81.11"
$ws.Range("E45").Value = "  -4.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.27"
$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "71.27"
$ws.Range("E47").Value = "  +5.25%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.97"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.37"
$ws.Range("E49").Value = "  +5.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("E50").Value = "  +7.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.17"
$ws.Range("E51").Value = "  +1.52%  "
